$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.234.96"
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = "'3.091.65"
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'523.55"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = "'136.31"
$ws.Range("E6").Value = '  -3.80%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'3.092.76"
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = '  +1.58%  '
$ws.Range("D13").Value = "'3.630.56"
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").Value = "'25.18"
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").Value = "'57.351.40"
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = "'3.094.12"
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = "'5.91"
$ws.Range("E19").Value = '  -2.88%  '
$ws.Range("D20").Value = "'12.32"
$ws.Range("E20").Value = '  -3.70%  '
$ws.Range("D21").Value = "'7.81"
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("D22").Value = "'346.82"
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = "'67.56"
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").Value = "'0.497"
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").Value = "'0.0₃0883"
$ws.Range("E28").Value = '  -3.83%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = "'7.32"
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").Value = "'1.86"
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = "'5.99"
$ws.Range("E32").Value = '  -7.73%  '
$ws.Range("D33").Value = "'20.64"
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("E34").Value = '  +6.00%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").Value = "'158.80"
$ws.Range("E36").Value = '  +1.82%  '
$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("D38").Value = "'25.70"
$ws.Range("E38").Value = '  -5.09%  '
$ws.Range("D39").Value = "'1.22"
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("E40").Value = '  +5.58%  '
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("E43").Value = '  +2.04%  '
$ws.Range("D44").Value = "'2.364.25"
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("D45").Value = "'36.47"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = "'0.0266"
$ws.Range("E47").Value = '  +2.20%  '
$ws.Range("D48").Value = "'0.965"
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("D49").Value = "'5.93"
$ws.Range("E49").Value = '  -1.51%  '
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("D51").Value = "'0.752"
$ws.Range("E51").Value = '  +2.47%  '
